$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row: add P1, Q1 with same style as O1 (bold/centered/bordered header style)
$ws.Range("P1").Value = 14
$ws.Range("Q1").Value = 15
$ws.Range("O1").Copy()
$ws.Range("P1:Q1").PasteSpecial(-4122)

# Data rows 2-25: swap I/K and M/O values, then append new P/Q columns = 2
for ($r = 2; $r -le 25; $r++) {
    $ws.Cells.Item($r, 9).Value = 2   # I -> 2
    $ws.Cells.Item($r, 11).Value = 1  # K -> 1
    $ws.Cells.Item($r, 13).Value = 2  # M -> 2
    $ws.Cells.Item($r, 15).Value = 1  # O -> 1
    $ws.Cells.Item($r, 16).Value = 2  # P -> 2
    $ws.Cells.Item($r, 17).Value = 2  # Q -> 2
}
